# Apply the commit's edits to the "Property" worksheet of Player.xlsx
#  1. Flip the "Save" flag (column E) OFF for the core player-stat rows (44-67)
#  2. Clean up the stray red/highlighted formatting that had been applied to the
#     two newly-added rows (76 "GameID", 77 "GateID") so they look like the
#     other plain rows (e.g. row 78), restoring the Text format on A76.
#  3. Leave the cursor/selection on H78, matching where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")
$ws.Activate()

# 1. Column E ("Save") TRUE -> FALSE for rows 44 through 67 (inclusive)
$ws.Range("E44:E67").Value = $false

# 2. Remove the highlighted style from rows 76 and 77, matching the
#    unformatted look of the surrounding rows (e.g. row 78).
$ws.Range("A76:J77").ClearFormats()

# A76 keeps a Text number format (like column A elsewhere, e.g. A2) even
# though the highlight is gone.
$ws.Range("A76").NumberFormat = $ws.Range("A2").NumberFormat

# 3. Update the active selection to H78.
$ws.Range("H78").Select()
